$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force the Price (D) and Volume(1h) (E) columns to be treated as text
# so that numeric-looking strings (e.g. '573.80') are not auto-converted
# into floating point numbers by Excel's COM value parser.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '69.727.39'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '3.564.56'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '573.80'
$ws.Range("E5").Value = '  -3.19%  '
$ws.Range("D6").Value = '186.24'
$ws.Range("E6").Value = '  -4.32%  '
$ws.Range("D7").Value = '3.563.98'
$ws.Range("E7").Value = '  -2.05%  '
$ws.Range("D8").Value = '0.620'
$ws.Range("E8").Value = '  -4.27%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '0.183'
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("D11").Value = '0.648'
$ws.Range("E11").Value = '  -3.96%  '
$ws.Range("D12").Value = '54.89'
$ws.Range("E12").Value = '  -6.16%  '
$ws.Range("D13").Value = '0.0000302'
$ws.Range("E13").Value = '  +3.15%  '
$ws.Range("D14").Value = '9.51'
$ws.Range("E14").Value = '  -4.32%  '
$ws.Range("D15").Value = '4.143.90'
$ws.Range("E15").Value = '  -2.54%  '
$ws.Range("D16").Value = '19.52'
$ws.Range("E16").Value = '  -1.95%  '
$ws.Range("D17").Value = '3.572.34'
$ws.Range("E17").Value = '  -2.37%  '
$ws.Range("D18").Value = '69.757.36'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("D19").Value = '12.45'
$ws.Range("E19").Value = '  -2.63%  '
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = '1.02'
$ws.Range("E21").Value = '  -4.00%  '
$ws.Range("D22").Value = '487.38'
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("D23").Value = '19.59'
$ws.Range("E23").Value = '  +2.88%  '
$ws.Range("D24").Value = '4.92'
$ws.Range("E24").Value = '  -7.89%  '
$ws.Range("D25").Value = '95.61'
$ws.Range("E25").Value = '  +4.66%  '
$ws.Range("D26").Value = '4.30'
$ws.Range("E26").Value = '  -4.55%  '
$ws.Range("D27").Value = '2.94'
$ws.Range("E27").Value = '  -7.16%  '
$ws.Range("D28").Value = '10.99'
$ws.Range("E28").Value = '  -6.15%  '
$ws.Range("D29").Value = '9.24'
$ws.Range("E29").Value = '  -4.04%  '
$ws.Range("D30").Value = '31.64'
$ws.Range("E30").Value = '  -3.70%  '
$ws.Range("D31").Value = '7.48'
$ws.Range("E31").Value = '  -5.58%  '
$ws.Range("D32").Value = '66.43'
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("D33").Value = '11.98'
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("D34").Value = '0.115'
$ws.Range("E34").Value = '  -5.06%  '
$ws.Range("D35").Value = '569.82'
$ws.Range("E35").Value = '  -9.74%  '
$ws.Range("D36").Value = '38.16'
$ws.Range("E36").Value = '  -6.35%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = '3.06'
$ws.Range("E38").Value = '  +6.24%  '
$ws.Range("D39").Value = '0.0₃0794'
$ws.Range("E39").Value = '  -4.13%  '
$ws.Range("D40").Value = '0.389'
$ws.Range("E40").Value = '  -5.89%  '
$ws.Range("D41").Value = '3.25'
$ws.Range("E41").Value = '  +12.00%  '
$ws.Range("D42").Value = '3.49'
$ws.Range("E42").Value = '  -2.94%  '
$ws.Range("E43").Value = '  -8.55%  '
$ws.Range("D44").Value = '3.255.67'
$ws.Range("E44").Value = '  -1.97%  '
$ws.Range("D45").Value = '2.98'
$ws.Range("E45").Value = '  -5.78%  '
$ws.Range("D46").Value = '3.47'
$ws.Range("E46").Value = '  +5.24%  '
$ws.Range("D47").Value = '0.0434'
$ws.Range("D48").Value = '9.63'
$ws.Range("E48").Value = '  +3.63%  '
$ws.Range("D49").Value = '0.135'
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").Value = '3.16'
$ws.Range("E51").Value = '  -4.94%  '

# Restore the original (default/Normal) cell style now that the text
# values are safely stored, so no stray number-format styling remains.
$dataRange.Style = "Normal"

